$d = $word.ActiveDocument

# Locate the paragraph that starts with "1. El juego debe arrancar..." and
# ends with "...cuando el jugador supera el actual." — this is the whole
# requirement #1 paragraph that needs to be struck through (paused).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*El juego debe arrancar con el nivel inicial*") {
        $r = $p.Range
        $r.Font.StrikeThrough = $true
        $r.Font.Size = 8
    }
}
